$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A56").Value = "2025-04-29 07:17:56"
$ws.Range("B56").Value = 151
